$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells are plain text in the source data (coin names, URLs,
# formatted prices and percentage strings). Force Text number format on
# each target cell before writing so numeric-looking values (e.g. "1.00",
# "615.91") are stored verbatim as text instead of being coerced to Number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.413.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.770.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.91'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.60'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.769.78'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.48'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.484'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.80'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.399.86'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.763.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.463.52'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.54'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '508.67'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.41'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.75%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.49'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.17'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +5.98%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.53'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.51%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.73%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.13%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.01'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.67%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.339'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '466.94'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.16%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.03'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +9.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.87'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.29'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.75%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.950.31'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.38'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.27'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.18%  '
